$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the tooltip for the "Python" skill (row 2) from "1+ years experience"
# to "2+ years experience" to reflect the added project experience.
$ws.Range("E2").Value = "2+ years experience"

# Move/restore the active selection to E2 (matches the saved state in the workbook).
$ws.Range("E2").Select()
